# "edições para trabalhar no caos"
#
# Re-create the header-row formatting + data-fix pass that was made to this
# sheet: recolor the header text to explicit black, replace the numeric
# placeholder in G2 with the text "NA", restore the (auto-fit) column
# widths, and leave the header row selected the way the sheet was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row (A1:J1): force an explicit black font color -----------
# (this is what introduces the 2nd <font>/<cellXfs> entry and the s="1"
# style index on every header cell)
$ws.Range("A1:J1").Font.Color = 0

# --- G2: the placeholder numeric 0 becomes the text value "NA" --------
$ws.Range("G2").Value = "NA"

# --- column widths: re-apply the (mac) auto-fit widths ------------------
$ws.Columns("A").ColumnWidth = 5.333333333333333
$ws.Columns("B").ColumnWidth = 8.5
$ws.Columns("C").ColumnWidth = 6
$ws.Columns("D").ColumnWidth = 3.5
$ws.Columns("E").ColumnWidth = 74.66666666666667
$ws.Columns("F").ColumnWidth = 3.5
$ws.Columns("G").ColumnWidth = 2.3333333333333335
$ws.Columns("H:I").ColumnWidth = 3.5
$ws.Columns("J").ColumnWidth = 21.5

# --- selection: sheet was left with the header row selected -----------
$ws.Range("A1:J1").Select()
